# Unify the conception of DataNode, DataTable, Entity.
# The worksheet previously named "Property1" is renamed to "DataNode".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"
